$p = $ppt.ActivePresentation

# --- Update the cached "datetimeFigureOut" date field text (10/22/2021 -> 10/25/2021) ---
# These placeholders live on the slide master, its first three layouts, and the notes master.
$newDate = "10/25/2021"

$master = $p.SlideMaster
$master.Shapes.Item("Date Placeholder 3").TextFrame.TextRange.Text = $newDate

for ($i = 1; $i -le 3; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    $layout.Shapes.Item("Date Placeholder 3").TextFrame.TextRange.Text = $newDate
}

try {
    $notesMaster = $p.NotesMaster
    $notesMaster.Shapes.Item("Date Placeholder 2").TextFrame.TextRange.Text = $newDate
} catch {
    Write-Host "NotesMaster date field could not be updated:" $_.Exception.Message
}

# --- Update slide titles: drop the "- sample 2" suffix ---
$p.Slides.Item(5).Shapes.Item(1).TextFrame.TextRange.Text = "Clustered column graph"
$p.Slides.Item(8).Shapes.Item(1).TextFrame.TextRange.Text = "Stacked Bar"
